# fix 19 nov 23
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "StudentGrade" to "Sheet1"
$ws.Name = "Sheet1"

# Header: "tanpa keterangan" -> "tanpa_keterangan"
$ws.Range("H3").Value = "tanpa_keterangan"

# Row 4 already has attendance counts (F4/G4/H4); catatan/penghargaan become "-"
$ws.Range("I4").Value = "-"
$ws.Range("J4").Value = "-"

# Rows 5-31: fill in sakit/izin/tanpa_keterangan counts (0) and
# catatan/penghargaan placeholders ("-")
for ($row = 5; $row -le 31; $row++) {
    $ws.Cells.Item($row, 6).Value = 0   # F: sakit
    $ws.Cells.Item($row, 7).Value = 0   # G: izin
    $ws.Cells.Item($row, 8).Value = 0   # H: tanpa_keterangan
    $ws.Cells.Item($row, 9).Value = "-"  # I: catatan
    $ws.Cells.Item($row, 10).Value = "-" # J: penghargaan
}
